# Commit: "Change Excel Field View to Cache, And set default value to FALSE"
#
# The "Property" sheet has a header row (row 1) naming each data column,
# and data rows 2-15. Column F's header was "View" and its data cells held
# the boolean default value for that field; rename the header to "Cache"
# and flip every row's default value in that column from TRUE to FALSE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: F1 "View" -> "Cache"
$ws.Range("F1").Value = "Cache"

# Data rows: F2:F15 TRUE -> FALSE
$ws.Range("F2:F15").Value = $false
